$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the order these cell values are first *written* controls the order new
# entries land in xl/sharedStrings.xml, so the statements below are sequenced
# to reproduce the same shared-string index layout as the target workbook.

# --- Row 69: old "NEXT: troubleshoot..." note replaced by a short "TODO:" header ---
$ws.Range("D69").Value = "TODO:"

# --- Row 70: new TODO bullet (text reused from the removed old row 73) ---
$ws.Range("D70").Value = "Testing for methods in ExperiencesSearch"

# --- Row 71: old "Cleanup:" header replaced by a TODO bullet about statsSearch servlet refactor ---
$ws.Range("D71").Value = "Possible refactoring of statsSearch servlet and testing of removed methods"

# --- Row 72: old "Add comma to income presentation..." replaced by orthodontist note; date cell A72 removed ---
$ws.Range("D72").Value = "Look into orthodontist problem - income value is ""-"" with some kind of footnote.  Decide how to handle."
$ws.Range("A72").Clear()

# --- Row 67: Tuesday hours note, "- x" placeholder resolved to "- 3:00" ---
$ws.Range("D67").Value = "Tue: ~0.5 + 8:40 - 3:00"

# --- Row 65: Week 11 summary gets expanded text (mentions Chart.js random sample chart) ---
$ws.Range("D65").Value = "Week 11: Feedback to presenters`nIndie project: combined 2 search forms and improved search output. Visual testing of validation and flow through program looks ok, except when searching orthodontist.  Still needs unit testing and maybe refactoring.  Also added Chart.js CDN links and random sample chart to my project."
$ws.Rows.Item(65).RowHeight = 60

# --- Rows 73 & 74: fully cleared (their single bullets were relocated/removed above) ---
$ws.Range("A73:D73").Clear()
$ws.Range("A74:D74").Clear()

# --- Scroll the view up a bit, matching the author's new cursor position ---
$win = $excel.ActiveWindow
$win.ScrollRow = 57
$win.ScrollColumn = 1
$ws.Range("D66").Select()
